# Auto-generated edit script: updates cached Leve profit-calculator values
# per the scheduled-runner refresh (H/I/J/K/L/M/N columns) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 90
$ws.Cells.Item(4, 9).Value = 90
$ws.Cells.Item(4, 11).Value = 90
$ws.Cells.Item(4, 13).Value = 24
$ws.Cells.Item(5, 8).Value = 441.58334
$ws.Cells.Item(5, 9).Value = 457.14285
$ws.Cells.Item(5, 11).Value = 457.14285
$ws.Cells.Item(5, 13).Value = -342.14285
$ws.Cells.Item(19, 8).Value = 1248.6666
$ws.Cells.Item(19, 10).Value = 1480.5
$ws.Cells.Item(19, 12).Value = 1480.5
$ws.Cells.Item(19, 14).Value = -1830.5
$ws.Cells.Item(33, 8).Value = 12500331
$ws.Cells.Item(33, 9).Value = 16666858
$ws.Cells.Item(33, 10).Value = 750
$ws.Cells.Item(33, 11).Value = 16666858
$ws.Cells.Item(33, 12).Value = 750
$ws.Cells.Item(33, 13).Value = -16666629
$ws.Cells.Item(33, 14).Value = -1208
$ws.Cells.Item(100, 8).Value = 1441.6923
$ws.Cells.Item(100, 9).Value = 1630.2222
$ws.Cells.Item(100, 10).Value = 1017.5
$ws.Cells.Item(100, 11).Value = 1630.2222
$ws.Cells.Item(100, 12).Value = 1017.5
$ws.Cells.Item(100, 13).Value = -1089.2222
$ws.Cells.Item(100, 14).Value = -2099.5
$ws.Cells.Item(106, 8).Value = 125003250
$ws.Cells.Item(106, 9).Value = 166669330
$ws.Cells.Item(106, 11).Value = 166669330
$ws.Cells.Item(106, 13).Value = -166668699
$ws.Cells.Item(113, 8).Value = 67972.39
$ws.Cells.Item(113, 10).Value = 15713.454
$ws.Cells.Item(113, 12).Value = 15713.454
$ws.Cells.Item(113, 14).Value = -22221.454
$ws.Cells.Item(125, 8).Value = 3445.5
$ws.Cells.Item(125, 9).Value = 2391
$ws.Cells.Item(125, 10).Value = 4500
$ws.Cells.Item(125, 11).Value = 21519
$ws.Cells.Item(125, 12).Value = 40500
$ws.Cells.Item(125, 13).Value = -19059
$ws.Cells.Item(125, 14).Value = -45420
$ws.Cells.Item(132, 8).Value = 9289213
$ws.Cells.Item(132, 9).Value = 11906367
$ws.Cells.Item(132, 10).Value = 129174.5
$ws.Cells.Item(132, 11).Value = 35719101
$ws.Cells.Item(132, 12).Value = 387523.5
$ws.Cells.Item(132, 13).Value = -35716571
$ws.Cells.Item(132, 14).Value = -392583.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3263.0107
$ws.Cells.Item(32, 9).Value = 1670.3292
$ws.Cells.Item(32, 11).Value = 1670.3292
$ws.Cells.Item(32, 13).Value = -1383.3292
$ws.Cells.Item(35, 8).Value = 1737.5
$ws.Cells.Item(35, 9).Value = 1737.5
$ws.Cells.Item(35, 11).Value = 1737.5
$ws.Cells.Item(35, 13).Value = -1331.5
$ws.Cells.Item(36, 8).Value = 27095.8
$ws.Cells.Item(36, 9).Value = 8869.75
$ws.Cells.Item(36, 10).Value = 100000
$ws.Cells.Item(36, 11).Value = 8869.75
$ws.Cells.Item(36, 12).Value = 100000
$ws.Cells.Item(36, 13).Value = -8523.75
$ws.Cells.Item(36, 14).Value = -100692
$ws.Cells.Item(61, 8).Value = 4220.7075
$ws.Cells.Item(61, 9).Value = 3257.8635
$ws.Cells.Item(61, 11).Value = 3257.8635
$ws.Cells.Item(61, 13).Value = -3045.8635
$ws.Cells.Item(110, 8).Value = 3846.2559
$ws.Cells.Item(110, 9).Value = 3589.1843
$ws.Cells.Item(110, 11).Value = 3589.1843
$ws.Cells.Item(110, 13).Value = -1544.1843
$ws.Cells.Item(136, 8).Value = 4220.7075
$ws.Cells.Item(136, 9).Value = 3257.8635
$ws.Cells.Item(136, 11).Value = 9773.5905
$ws.Cells.Item(136, 13).Value = -7223.5905

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(92, 8).Value = 40000
$ws.Cells.Item(92, 10).Value = 40000
$ws.Cells.Item(92, 12).Value = 40000
$ws.Cells.Item(92, 14).Value = -44992
$ws.Cells.Item(97, 8).Value = 12790.3
$ws.Cells.Item(97, 9).Value = 8491
$ws.Cells.Item(97, 11).Value = 8491
$ws.Cells.Item(97, 13).Value = -7500
$ws.Cells.Item(134, 8).Value = 1845.4468
$ws.Cells.Item(134, 9).Value = 1812.2046
$ws.Cells.Item(134, 10).Value = 2333
$ws.Cells.Item(134, 11).Value = 5436.6138
$ws.Cells.Item(134, 12).Value = 6999
$ws.Cells.Item(134, 13).Value = -2901.6138
$ws.Cells.Item(134, 14).Value = -12069

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 218.6875
$ws.Cells.Item(7, 9).Value = 235.64285
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 235.64285
$ws.Cells.Item(7, 12).Value = 100
$ws.Cells.Item(7, 13).Value = -122.64285
$ws.Cells.Item(7, 14).Value = -326
$ws.Cells.Item(64, 8).Value = 80000
$ws.Cells.Item(64, 10).Value = 80000
$ws.Cells.Item(64, 12).Value = 80000
$ws.Cells.Item(64, 14).Value = -80496
$ws.Cells.Item(67, 8).Value = 80000
$ws.Cells.Item(67, 10).Value = 80000
$ws.Cells.Item(67, 12).Value = 80000
$ws.Cells.Item(67, 14).Value = -81716
$ws.Cells.Item(68, 8).Value = 77333.2
$ws.Cells.Item(68, 10).Value = 100334.5
$ws.Cells.Item(68, 12).Value = 100334.5
$ws.Cells.Item(68, 14).Value = -101832.5
$ws.Cells.Item(71, 8).Value = 77333.2
$ws.Cells.Item(71, 10).Value = 100334.5
$ws.Cells.Item(71, 12).Value = 301003.5
$ws.Cells.Item(71, 14).Value = -308491.5
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 14).Value = ""
$ws.Cells.Item(107, 8).Value = 1206.931
$ws.Cells.Item(107, 10).Value = 1229.1111
$ws.Cells.Item(107, 12).Value = 1229.1111
$ws.Cells.Item(107, 14).Value = -5069.1111
$ws.Cells.Item(121, 8).Value = 49666.668
$ws.Cells.Item(121, 10).Value = 49666.668
$ws.Cells.Item(121, 12).Value = 49666.668
$ws.Cells.Item(121, 14).Value = -52286.668
$ws.Cells.Item(134, 8).Value = 15600.149
$ws.Cells.Item(134, 9).Value = 13910.019
$ws.Cells.Item(134, 10).Value = 23346.584
$ws.Cells.Item(134, 11).Value = 41730.057
$ws.Cells.Item(134, 12).Value = 70039.75199999999
$ws.Cells.Item(134, 13).Value = -39195.057
$ws.Cells.Item(134, 14).Value = -75109.75199999999
$ws.Cells.Item(135, 8).Value = 48738.8
$ws.Cells.Item(135, 10).Value = 48738.8
$ws.Cells.Item(135, 12).Value = 48738.8
$ws.Cells.Item(135, 14).Value = -58878.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 359.375
$ws.Cells.Item(17, 9).Value = 469.33334
$ws.Cells.Item(17, 11).Value = 1408.00002
$ws.Cells.Item(17, 13).Value = -1239.00002
$ws.Cells.Item(117, 8).Value = 1259
$ws.Cells.Item(117, 9).Value = 825
$ws.Cells.Item(117, 10).Value = 1548.3334
$ws.Cells.Item(117, 11).Value = 2475
$ws.Cells.Item(117, 12).Value = 4645.0002
$ws.Cells.Item(117, 13).Value = 967
$ws.Cells.Item(117, 14).Value = -11529.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 25642100
$ws.Cells.Item(102, 9).Value = 1116.9117
$ws.Cells.Item(102, 10).Value = 200000800
$ws.Cells.Item(102, 11).Value = 1116.9117
$ws.Cells.Item(102, 12).Value = 200000800
$ws.Cells.Item(102, 13).Value = 505.0882999999999
$ws.Cells.Item(102, 14).Value = -200004044
$ws.Cells.Item(113, 8).Value = 2068.5386
$ws.Cells.Item(113, 10).Value = 1970.5
$ws.Cells.Item(113, 12).Value = 1970.5
$ws.Cells.Item(113, 14).Value = -6310.5
$ws.Cells.Item(122, 8).Value = 3015.95
$ws.Cells.Item(122, 9).Value = 2510.1538
$ws.Cells.Item(122, 10).Value = 3955.2856
$ws.Cells.Item(122, 11).Value = 7530.4614
$ws.Cells.Item(122, 12).Value = 11865.8568
$ws.Cells.Item(122, 13).Value = -5080.4614
$ws.Cells.Item(122, 14).Value = -16765.8568
$ws.Cells.Item(124, 8).Value = 29999
$ws.Cells.Item(124, 10).Value = 29999
$ws.Cells.Item(124, 12).Value = 29999
$ws.Cells.Item(124, 14).Value = -39819
$ws.Cells.Item(126, 8).Value = 3780.5264
$ws.Cells.Item(126, 9).Value = 3662.0667
$ws.Cells.Item(126, 10).Value = 4224.75
$ws.Cells.Item(126, 11).Value = 10986.2001
$ws.Cells.Item(126, 12).Value = 12674.25
$ws.Cells.Item(126, 13).Value = -8516.2001
$ws.Cells.Item(126, 14).Value = -17614.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4901.864
$ws.Cells.Item(7, 9).Value = 4209.125
$ws.Cells.Item(7, 11).Value = 4209.125
$ws.Cells.Item(7, 13).Value = -4097.125
$ws.Cells.Item(31, 8).Value = 3615.55
$ws.Cells.Item(31, 10).Value = 7004.6
$ws.Cells.Item(31, 12).Value = 7004.6
$ws.Cells.Item(31, 14).Value = -7500.6
$ws.Cells.Item(40, 8).Value = 4006.2666
$ws.Cells.Item(40, 9).Value = 3596.0264
$ws.Cells.Item(40, 10).Value = 6233.2856
$ws.Cells.Item(40, 11).Value = 3596.0264
$ws.Cells.Item(40, 12).Value = 6233.2856
$ws.Cells.Item(40, 13).Value = -3460.0264
$ws.Cells.Item(40, 14).Value = -6505.2856
$ws.Cells.Item(68, 8).Value = 4266.222
$ws.Cells.Item(68, 9).Value = 4670.857
$ws.Cells.Item(68, 10).Value = 2850
$ws.Cells.Item(68, 11).Value = 4670.857
$ws.Cells.Item(68, 12).Value = 2850
$ws.Cells.Item(68, 13).Value = -3921.857
$ws.Cells.Item(68, 14).Value = -4348
$ws.Cells.Item(71, 8).Value = 4266.222
$ws.Cells.Item(71, 9).Value = 4670.857
$ws.Cells.Item(71, 10).Value = 2850
$ws.Cells.Item(71, 11).Value = 23354.285
$ws.Cells.Item(71, 12).Value = 14250
$ws.Cells.Item(71, 13).Value = -19610.285
$ws.Cells.Item(71, 14).Value = -21738
$ws.Cells.Item(100, 8).Value = 1983.6666
$ws.Cells.Item(100, 10).Value = 2499.5
$ws.Cells.Item(100, 12).Value = 2499.5
$ws.Cells.Item(100, 14).Value = -3581.5
$ws.Cells.Item(122, 8).Value = 440455.8
$ws.Cells.Item(122, 9).Value = 629711.6
$ws.Cells.Item(122, 10).Value = 7871.143
$ws.Cells.Item(122, 11).Value = 1889134.8
$ws.Cells.Item(122, 12).Value = 23613.429
$ws.Cells.Item(122, 13).Value = -1886684.8
$ws.Cells.Item(122, 14).Value = -28513.429
$ws.Cells.Item(126, 8).Value = 4901.864
$ws.Cells.Item(126, 9).Value = 4209.125
$ws.Cells.Item(126, 11).Value = 12627.375
$ws.Cells.Item(126, 13).Value = -10157.375
$ws.Cells.Item(132, 8).Value = 3394.2334
$ws.Cells.Item(132, 9).Value = 3070.2693
$ws.Cells.Item(132, 10).Value = 5500
$ws.Cells.Item(132, 11).Value = 9210.8079
$ws.Cells.Item(132, 12).Value = 16500
$ws.Cells.Item(132, 13).Value = -6680.8079
$ws.Cells.Item(132, 14).Value = -21560
$ws.Cells.Item(133, 8).Value = 55825
$ws.Cells.Item(133, 9).Value = 45000
$ws.Cells.Item(133, 10).Value = 57990
$ws.Cells.Item(133, 11).Value = 45000
$ws.Cells.Item(133, 12).Value = 57990
$ws.Cells.Item(133, 13).Value = -42470
$ws.Cells.Item(133, 14).Value = -63050
$ws.Cells.Item(136, 8).Value = 2306.55
$ws.Cells.Item(136, 9).Value = 1971.7241
$ws.Cells.Item(136, 11).Value = 5915.1723
$ws.Cells.Item(136, 13).Value = -3365.1723

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 8000
$ws.Cells.Item(18, 9).Value = 8000
$ws.Cells.Item(18, 10).Value = 8000
$ws.Cells.Item(18, 11).Value = 8000
$ws.Cells.Item(18, 12).Value = 8000
$ws.Cells.Item(18, 13).Value = -7827
$ws.Cells.Item(18, 14).Value = -8346
$ws.Cells.Item(107, 8).Value = 825.5599999999999
$ws.Cells.Item(107, 9).Value = 857.25
$ws.Cells.Item(107, 10).Value = 698.8
$ws.Cells.Item(107, 11).Value = 2571.75
$ws.Cells.Item(107, 12).Value = 2096.4
$ws.Cells.Item(107, 13).Value = -651.75
$ws.Cells.Item(107, 14).Value = -5936.4
$ws.Cells.Item(122, 8).Value = 2819.5334
$ws.Cells.Item(122, 10).Value = 3036.875
$ws.Cells.Item(122, 12).Value = 9110.625
$ws.Cells.Item(122, 14).Value = -14010.625

